# Rotate the species-specific data among rows 6, 7 and 8:
#   new row6 = old row7 data
#   new row7 = old row8 data
#   new row8 = old row6 data
# Only columns A, B, E, F, G, H, Q, R differ between the three rows; the
# remaining columns (C, D, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG,
# AT, AW, AX, AY) are identical across rows 6-8 and therefore unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Capture the original values for rows 6, 7 and 8 before overwriting anything.
# Note: use .Value2 (not .Value) to read/write - .Value is unreliable in this
# COM-interop runtime.
$orig = @{}
foreach ($r in 6, 7, 8) {
    $orig[$r] = @{}
    foreach ($col in $cols) {
        $orig[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# row6 <- row7, row7 <- row8, row8 <- row6
$mapping = @{ 6 = 7; 7 = 8; 8 = 6 }

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $orig[$srcRow][$col]
    }
}
